$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '25.904.82'
$ws.Range("E2").Value = '  +0.22%  '

# Row 3
$ws.Range("D3").Value = '1.639.50'
$ws.Range("E3").Value = '  +0.16%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("E4").Value = '  +0.32%  '

# Row 5
$ws.Range("E5").Value = '  -0.32%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5043'
$ws.Range("E6").Value = '  +0.43%  '

# Row 8
$ws.Range("E8").Value = '  +0.11%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06399'
$ws.Range("E9").Value = '  -0.29%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.66'
$ws.Range("E10").Value = '  +0.28%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07787'
$ws.Range("E11").Value = '  +1.31%  '

# Row 12
$ws.Range("E12").Value = '  +0.94%  '

# Row 13
$ws.Range("D13").Value = '1.649.06'
$ws.Range("E13").Value = '  +0.74%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5438'
$ws.Range("E14").Value = '  -0.45%  '

# Row 15
$ws.Range("D15").Value = '0.0₅7876'
$ws.Range("E15").Value = '  -0.62%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.07'
$ws.Range("E16").Value = '  +2.40%  '

# Row 17
$ws.Range("D17").Value = '25.956.73'
$ws.Range("E17").Value = '  +0.36%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.004'
$ws.Range("E18").Value = '  +0.12%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '197.73'
$ws.Range("E19").Value = '  -2.76%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.392'
$ws.Range("E20").Value = '  +2.07%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.963'
$ws.Range("E21").Value = '  +0.07%  '

# Row 22
$ws.Range("E22").Value = '  +0.52%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.007'
$ws.Range("E23").Value = '  +0.27%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.870'
$ws.Range("E24").Value = '  -3.39%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '140.24'
$ws.Range("E25").Value = '  -0.84%  '

# Row 26
$ws.Range("E26").Value = '  -0.46%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.866'
$ws.Range("E27").Value = '  +2.35%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.71'
$ws.Range("E28").Value = '  +0.26%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.238'

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05015'
$ws.Range("E30").Value = '  +1.08%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.261'
$ws.Range("E31").Value = '  -0.28%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.194'
$ws.Range("E32").Value = '  +0.11%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.538'
$ws.Range("E33").Value = '  +0.52%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.363'
$ws.Range("E34").Value = '  +0.46%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.8940'
$ws.Range("E35").Value = '  +0.02%  '

# Row 36
$ws.Range("E36").Value = '  -1.01%  '

# Row 37
$ws.Range("D37").Value = '1.136.58'
$ws.Range("E37").Value = '  -3.45%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5526'

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01555'
$ws.Range("E39").Value = '  -0.17%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.005'
$ws.Range("E40").Value = '  +0.23%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.700'
$ws.Range("E41").Value = '  +1.09%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8148'
$ws.Range("E42").Value = '  +1.31%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '99.44'
$ws.Range("E43").Value = '  -0.12%  '

# Row 44
$ws.Range("E44").Value = '  +10.54%  '

# Row 45
$ws.Range("D45").Value = '1.777.46'
$ws.Range("E45").Value = '  +0.24%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4530'
$ws.Range("E46").Value = '  +0.40%  '

# Row 47
$ws.Range("B47").Value = 'Frax'
$ws.Range("C47").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.006'
$ws.Range("E47").Value = '  +0.01%  '

# Row 48
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '55.20'
$ws.Range("E48").Value = '  +0.67%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05082'
$ws.Range("E49").Value = '  +0.70%  '

# Row 50
$ws.Range("B50").Value = 'USDD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.005'
$ws.Range("E50").Value = '  +0.39%  '

# Row 51
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.09531'
$ws.Range("E51").Value = '  +2.94%  '
